$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# --- Update the capacity-market rerun inputs (rows 28-30) ---
# fix_demand_to_representative_year: TRUE -> FALSE
$ws.Range("B28").Value = $false
# fix_profiles_to_representative_year: TRUE -> FALSE
$ws.Range("B29").Value = $false
# Representative year: 2015 -> 2004
$ws.Range("B30").Value = 2004

# Highlight the now-active inputs with the yellow fill used elsewhere
# in the sheet for "active" settings (style index 7 / RGB FFFFFF00).
$ws.Range("B28").Interior.Color = 65535
$ws.Range("B29").Interior.Color = 65535
$ws.Range("B30").Interior.Color = 65535

# --- Restore the view state (scrolled down a bit further, new selection) ---
$ws.Activate()
$ws.Range("C30").Select()
